# Initial check-in of translations changes.
#
# The "display.*" translation-key column headers on the survey and settings
# sheets are renamed to their new ".text" suffixed keys, and the workbook is
# left with "settings" as the active/selected sheet (instead of "survey").

$wb = $excel.ActiveWorkbook

# --- survey sheet: rename translation-key headers (row 1) ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("E1").Value = "display.hint.text"
$survey.Range("D1").Value = "display.prompt.text"

# --- settings sheet: rename translation-key header (row 1) ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"

# --- restore/update the selected cell on each sheet ---
[void]$survey.Range("D2").Select()
[void]$settings.Range("C2").Select()

# --- "settings" becomes the active sheet (tab selected) on reopen ---
[void]$settings.Activate()
